$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting existing C:E to D:F
$ws.Columns.Item(3).Insert()

# Match the (approximate) width of the new column C to column B
$ws.Columns.Item(3).ColumnWidth = 31.5

# Apply the bold/shaded header style (style index 2, same as the rest of row 1) to C1
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Apply the value-row style (style index 1, same as B2) to C2
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Set the value first (so "member for" lands before the header text in the
# shared string table), then the header text
$ws.Range("C2").Value = "member for"
$ws.Range("C1").Value = "Account_Membership_Term"

# Update the selection to match the target workbook view
$ws.Range("C1").Select()
